$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers (B1:Z1), using the same style as A1 (bold + centered header) ---
$headers = @(
  "nea-ukcharity.bsky.social",
  "caneurope.bsky.social",
  "wateraid.bsky.social",
  "migrantsrights.bsky.social",
  "friends-earth.bsky.social",
  "samcardwell44.bsky.social",
  "greenpeace.eu",
  "wwfeu.bsky.social",
  "powertochange.org.uk",
  "thegreenregister.bsky.social",
  "endfuelpoverty.bsky.social",
  "commenergyengland.bsky.social",
  "extinctionrebellion.uk",
  "wwtworldwide.bsky.social",
  "bristolgreenparty.bsky.social",
  "warmthiswinter.bsky.social",
  "jrct.bsky.social",
  "ssencommunity.bsky.social",
  "localtrust.bsky.social",
  "wiltscouncil.bsky.social",
  "nationalgrid.bsky.social",
  "ofgem.bsky.social",
  "barnsleycouncil.bsky.social",
  "northsomersetc.bsky.social",
  "citizensadvice.bsky.social"
)

# --- Row 2 values (B2:Z2), aligned with the headers above ---
$values = @(293, 87, 17, 82, 0, 2, 372, 54, 108, 7, 102, 42, 72, 0, 650, 0, 32, 17, 81, 0, 0, 59, 0, 0, 26)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 2
    $headerCell = $ws.Cells.Item(1, $col)
    $headerCell.Value = $headers[$i]
    $headerCell.Font.Bold = $true
    $headerCell.HorizontalAlignment = -4108

    $ws.Cells.Item(2, $col).Value = $values[$i]
}

# --- A2: update date text, keeping it a plain (unstyled) text cell ---
$dateCell = $ws.Range("A2")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-01-21"
$dateCell.ClearFormats()

Write-Output "Done"
